# Apply the edits described by the diff:
#  - Slide 3 ("Törtenete"): last content bullet "...legújabb a 12-es"
#    becomes "...legújabb " / "a 11-es" (two runs, same visible text except 12->11).
#  - Slide 6 ("Hátrányok"): title becomes "Hátrányok (iOS-hez képest)" with a
#    trailing blank line; content placeholder is replaced with three new
#    bullets plus two trailing blank paragraphs.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 3: "Folyamatos fejlesztés, legújabb a 12-es" -> "... a 11-es"
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$body3 = $s3.Shapes.Item(2).TextFrame.TextRange
$paraCount3 = $body3.Paragraphs().Count
$lastPara = $body3.Paragraphs($paraCount3)

$full = $lastPara.Text
$idx = $full.IndexOf("a 12-es")
if ($idx -ge 0) {
    $target = $lastPara.Characters($idx + 1, 7)
    $target.Text = "a 11-es"
}

# ---------------------------------------------------------------------
# Slide 6: title "Hátrányok" -> "Hátrányok (iOS-hez képest)" + blank line
# ---------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$title6 = $s6.Shapes.Item(1).TextFrame.TextRange

$title6.InsertAfter(" (iOS-hez képest)")

$full = $title6.Text
$parenIdx = $full.IndexOf("(")
$prefix = $title6.Characters(1, $parenIdx + 1)
$prefix.Text = "Hátrányok ("

$full = $title6.Text
$iosIdx = $full.IndexOf("iOS-hez")
$iosRun = $title6.Characters($iosIdx + 1, 7)
$iosRun.Text = "iOS-hez"

# Trailing blank line after the parenthesis text.
$vtab = [char]11
$title6.InsertAfter($vtab)

# ---------------------------------------------------------------------
# Slide 6: content placeholder -> three bullets + two blank paragraphs
# ---------------------------------------------------------------------
$content6 = $s6.Shapes.Item(2).TextFrame.TextRange
$cr = [char]13
$newBody = "Nem elég tiszta a megjelenés " + $cr + "Nagyobb a hibalehetőség" + $cr + "Könnyebben feltörhető" + $cr + $cr
$content6.Text = $newBody
